$wb = $excel.ActiveWorkbook

# The same edits apply to both the "展览" (sheet 1) and "全部类型" (sheet 4)
# worksheets, which mirror each other's data.
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # --- Row 2 ---
    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = "2024.01.27"
    $ws.Range("B2").Style = "Normal"
    $ws.Range("C2").Value = "南宁·AP动漫游戏嘉年华"
    $ws.Range("F2").Value = 2012
    $ws.Range("I2").Value = "https://show.bilibili.com/platform/detail.html?id=79764&msource=Msearch_colligation"

    # --- Row 3 ---
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = "2024.01.27"
    $ws.Range("B3").Style = "Normal"
    $ws.Range("C3").Value = "南宁·第一届异次元动漫嘉年华"
    $ws.Range("F3").Value = 601
    $ws.Range("I3").Value = "https://show.bilibili.com/platform/detail.html?id=78089&msource=Msearch_colligation"

    # --- Row 4 ---
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024.02.01"
    $ws.Range("B4").Style = "Normal"
    $ws.Range("C4").Value = "南宁·桂南动漫游戏嘉年华"
    $ws.Range("F4").Value = 1360
    $ws.Range("H4").Value = $true
    $ws.Range("I4").Value = "https://show.bilibili.com/platform/detail.html?id=79354&msource=Msearch_colligation"

    # --- Row 5 ---
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2024.02.15"
    $ws.Range("B5").Style = "Normal"
    $ws.Range("C5").Value = "南宁·2024良牙动漫冬季盛典（冬典）"
    $ws.Range("F5").Value = 6740
    $ws.Range("I5").Value = "https://show.bilibili.com/platform/detail.html?id=77938&msource=Msearch_colligation"

    # --- Row 6 ---
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2024.03.09"
    $ws.Range("B6").Style = "Normal"
    $ws.Range("C6").Value = "南宁·第五届小蜜蜂动漫嘉年华"
    $ws.Range("F6").Value = 167
    $ws.Range("I6").Value = "https://show.bilibili.com/platform/detail.html?id=79051&msource=Msearch_colligation"

    # --- Row 7 ---
    $ws.Range("B7").NumberFormat = "@"
    $ws.Range("B7").Value = "2024.03.16"
    $ws.Range("B7").Style = "Normal"
    $ws.Range("C7").Value = "南宁·草莓动漫节"
    $ws.Range("F7").Value = 57
    $ws.Range("I7").Value = "https://show.bilibili.com/platform/detail.html?id=80943&msource=Msearch_colligation"
}
